# Auto-generated Excel COM-interop script
# Applies numeric cell updates/insertions/deletions across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1013.5
$ws.Range("J8").Value = 2888.6667
$ws.Range("L8").Value = 8666.000100000001
$ws.Range("N8").Value = -8944.000100000001
$ws.Range("H15").Value = 1026.7059
$ws.Range("I15").Value = 1026.7059
$ws.Range("K15").Value = 3080.1177
$ws.Range("M15").Value = -2911.1177
$ws.Range("H98").Value = 2866
$ws.Range("I98").Value = 1799
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 1799
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -301
$ws.Range("N98").Value = -7996
$ws.Range("H122").Value = 2866
$ws.Range("I122").Value = 1799
$ws.Range("K122").Value = 5397
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2947
$ws.Range("N122").Value = -19900
$ws.Range("H137").Value = 3271
$ws.Range("I137").Value = 3061.077
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 9183.231
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -6633.231
$ws.Range("N137").Value = -23100
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1473.7142
$ws.Range("I32").Value = 1473.7142
$ws.Range("K32").Value = 1473.7142
$ws.Range("M32").Value = -1186.7142
$ws.Range("H61").Value = 3523.111
$ws.Range("J61").Value = 4317.4614
$ws.Range("L61").Value = 4317.4614
$ws.Range("N61").Value = -4741.4614
$ws.Range("H74").Value = 992.5
$ws.Range("I74").Value = 992.5
$ws.Range("K74").Value = 992.5
$ws.Range("M74").Value = -118.5
$ws.Range("H77").Value = 992.5
$ws.Range("I77").Value = 992.5
$ws.Range("K77").Value = 4962.5
$ws.Range("M77").Value = -594.5
$ws.Range("H136").Value = 3523.111
$ws.Range("J136").Value = 4317.4614
$ws.Range("L136").Value = 12952.3842
$ws.Range("N136").Value = -18052.3842
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 563.8333
$ws.Range("I64").Value = 586.6
$ws.Range("K64").Value = 586.6
$ws.Range("M64").Value = -361.6
$ws.Range("H67").Value = 563.8333
$ws.Range("I67").Value = 586.6
$ws.Range("K67").Value = 586.6
$ws.Range("M67").Value = 193.4
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 1100
$ws.Range("K16").Value = 1100
$ws.Range("M16").Value = -813
$ws.Range("H50").Value = 25996.666
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 33995
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 33995
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = -35245
$ws.Range("H59").Value = 30371.25
$ws.Range("I59").Value = 21750
$ws.Range("J59").Value = 38992.5
$ws.Range("K59").Value = 21750
$ws.Range("L59").Value = 38992.5
$ws.Range("M59").Value = -20605
$ws.Range("N59").Value = -41282.5
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 782.1
$ws.Range("I8").Value = 782.1
$ws.Range("K8").Value = 2346.3
$ws.Range("M8").Value = -2207.3
$ws.Range("H12").Value = 166.66667
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 249.5
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 748.5
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -1094.5
$ws.Range("H26").Value = 999
$ws.Range("I26").Value = 999
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 2997
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -2709
$ws.Range("N26").Value = ""
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("H93").Value = 2999.5
$ws.Range("J93").Value = 2999.5
$ws.Range("L93").Value = 8998.5
$ws.Range("N93").Value = -12742.5
$ws.Range("H109").Value = 4791.8184
$ws.Range("I109").Value = 2710
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 8130
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -7090
$ws.Range("N109").Value = -17080
$ws.Range("H131").Value = 1599.6666
$ws.Range("J131").Value = 2181.8
$ws.Range("L131").Value = 6545.400000000001
$ws.Range("N131").Value = -16625.4
$ws.Range("H140").Value = 2511.9285
$ws.Range("I140").Value = 2511.9285
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 7535.7855
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2355.7855
$ws.Range("N140").Value = ""
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5555
$ws.Range("I70").Value = 5555
$ws.Range("K70").Value = 5555
$ws.Range("M70").Value = -5285
$ws.Range("H73").Value = 5555
$ws.Range("I73").Value = 5555
$ws.Range("K73").Value = 5555
$ws.Range("M73").Value = -4619
$ws.Range("H122").Value = 1641.4286
$ws.Range("J122").Value = 1497.5
$ws.Range("L122").Value = 4492.5
$ws.Range("N122").Value = -9392.5
$ws.Range("H132").Value = 3630.2727
$ws.Range("I132").Value = 3922.6667
$ws.Range("J132").Value = 3279.4
$ws.Range("K132").Value = 11768.0001
$ws.Range("L132").Value = 9838.200000000001
$ws.Range("M132").Value = -9238.000100000001
$ws.Range("N132").Value = -14898.2
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 3500003
$ws.Range("I13").Value = 3500003
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 3500003
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -3499863
$ws.Range("N13").Value = ""
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 833.3333
$ws.Range("I126").Value = 833.3333
$ws.Range("K126").Value = 2499.9999
$ws.Range("M126").Value = -29.9998999999998
$ws.Range("H132").Value = 1908.9
$ws.Range("I132").Value = 613
$ws.Range("K132").Value = 1839
$ws.Range("M132").Value = 691

Write-Host "Applied 160 cell updates"
